# IMPROVE: Schedule and jobs
# Update the error description and the date of the last activity report
# on the "CurrentCurrencyTrades" schedule/job sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (error) / Column G (date) on row 4 hold the job's last error text
# and timestamp. Replace them with the new message/timestamp.
$ws.Range("F4").Value = "Was not able to save the task of exchange rates."
$ws.Range("G4").Value = "2022-09-26 12:39:49"
